# Apply the cryptos list refresh described in the commit:
# "Updated cryptos list on Mon Sep 11 07:35:36 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like a plain number (e.g. "212.59").
# Force those specific cells to Text format first so Excel keeps them as strings
# instead of silently re-typing them as numeric values.
$textFormatCells = @(
    'D5',
    'D7',
    'D8',
    'D10',
    'D14',
    'D17',
    'D20',
    'D22',
    'D25',
    'D32',
    'D39',
    'D41',
    'D47',
    'D50',
    'D51'
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = '@'
}

# Updated price (column D) and 1h volume change (column E) values
$ws.Range('D2').Value = '25.945.72'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.620.61'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '212.59'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = '0.489'
$ws.Range('E7').Value = '  -3.33%  '
$ws.Range('D8').Value = '0.0621'
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').Value = '18.31'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '1.846.57'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').Value = '1.611.63'
$ws.Range('E13').Value = '  -3.18%  '
$ws.Range('D14').Value = '4.14'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').Value = '25.965.30'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = '61.68'
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '191.79'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D22').Value = '9.53'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '144.41'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -4.00%  '
$ws.Range('E28').Value = '  -1.97%  '
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('D32').Value = '3.11'
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range('D36').Value = '1.129.31'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('E37').Value = '  -3.96%  '
$ws.Range('E38').Value = '  -1.89%  '
$ws.Range('D39').Value = '0.518'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').Value = '97.82'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('D42').Value = '1.757.56'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('E43').Value = '  -4.02%  '
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('D47').Value = '54.04'
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').Value = '7.49'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  +0.10%  '
